# "new data for prediction" — refresh the rolling forecast window on the
# Maywood-Garfield meal-count sheet:
#   - clear a handful of placeholder/forecast cells back to blank
#   - update a few real numbers that came in
#   - drop the now-stale tail of the rolling window (11 rows)
#   - keep autoFilter / defined name / dimension in sync with the new extent
#   - move the selection to where the user was working

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- clear stale placeholder numbers back to blank (keeps cell styles) ----
$ws.Range("C229:E230").ClearContents()
$ws.Range("C236:E237").ClearContents()
$ws.Range("C247:E248").ClearContents()
$ws.Range("C254:E255").ClearContents()

# ---- real numbers landed for a few days ----
$ws.Range("C253").Value = 66
$ws.Range("D253").Value = 80
$ws.Range("C256").Value = 68
$ws.Range("D256").Value = 75
$ws.Range("C257").Value = 75
$ws.Range("D257").Value = 78

# ---- the rolling window's tail (rows 258:268) is now stale; drop it so ----
# ---- everything below shifts up, shedding 11 rows off the bottom        ----
$ws.Rows("258:268").Delete()

# ---- keep the autofilter + the hidden _FilterDatabase name in sync ----
$ws.AutoFilterMode = $false
$ws.Range("B1:B996").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Maywood-Garfield'!`$B`$1:`$B`$996"
    }
}

# ---- move the selection to where work continued ----
$ws.Range("A248:XFD248").Select()
